$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, [string]$text) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-TextValue "D2" "69.518.51"
Set-TextValue "E2" "  +2.99%  "
Set-TextValue "D3" "3.385.83"
Set-TextValue "E3" "  +4.67%  "
Set-TextValue "E4" "  +0.02%  "
Set-TextValue "D5" "191.87"
Set-TextValue "E5" "  +4.21%  "
Set-TextValue "D6" "593.12"
Set-TextValue "E6" "  +2.40%  "
Set-TextValue "E7" "  +0.02%  "
Set-TextValue "D8" "0.607"
Set-TextValue "E8" "  -0.58%  "
Set-TextValue "E9" "  +3.25%  "
Set-TextValue "D10" "6.78"
Set-TextValue "E10" "  +3.07%  "
Set-TextValue "E11" "  +2.30%  "
Set-TextValue "D12" "3.976.88"
Set-TextValue "E12" "  +4.91%  "
Set-TextValue "E13" "  +1.16%  "
Set-TextValue "D14" "28.73"
Set-TextValue "E14" "  +4.09%  "
Set-TextValue "D15" "69.568.62"
Set-TextValue "E15" "  +3.01%  "
Set-TextValue "D16" "0.0000172"
Set-TextValue "E16" "  +2.23%  "
Set-TextValue "D17" "3.400.16"
Set-TextValue "E17" "  +5.28%  "
Set-TextValue "D18" "450.38"
Set-TextValue "E18" "  +14.20%  "
Set-TextValue "E19" "  +1.72%  "
Set-TextValue "D20" "13.84"
Set-TextValue "E20" "  +2.77%  "
Set-TextValue "E21" "  +3.69%  "
Set-TextValue "E22" "  +5.41%  "
Set-TextValue "E23" "  +0.16%  "
Set-TextValue "D24" "3.523.62"
Set-TextValue "E24" "  +4.61%  "
Set-TextValue "E25" "  +4.51%  "
Set-TextValue "E26" "  +1.52%  "
Set-TextValue "E27" "  +1.60%  "
Set-TextValue "D28" "9.50"
Set-TextValue "E28" "  -1.18%  "
Set-TextValue "E29" "  +0.21%  "
Set-TextValue "E30" "  +1.82%  "
Set-TextValue "D31" "23.45"
Set-TextValue "E31" "  +3.91%  "
Set-TextValue "D32" "5.67"
Set-TextValue "E32" "  +2.05%  "
Set-TextValue "D33" "1.29"
Set-TextValue "E33" "  +3.47%  "
Set-TextValue "D34" "7.00"
Set-TextValue "E34" "  +0.11%  "
Set-TextValue "E35" "  +0.01%  "
Set-TextValue "E36" "  +5.74%  "
Set-TextValue "D37" "165.33"
Set-TextValue "E37" "  +2.45%  "
Set-TextValue "E38" "  +3.89%  "
Set-TextValue "D39" "27.56"
Set-TextValue "E39" "  +4.32%  "
Set-TextValue "D40" "0.819"
Set-TextValue "E40" "  +2.01%  "
Set-TextValue "D41" "4.61"
Set-TextValue "E41" "  +1.37%  "
Set-TextValue "D42" "6.57"
Set-TextValue "E42" "  +1.57%  "
Set-TextValue "D43" "2.758.80"
Set-TextValue "E43" "  +5.70%  "
Set-TextValue "E44" "  +3.01%  "
Set-TextValue "D45" "25.66"
Set-TextValue "E45" "  +3.65%  "
Set-TextValue "D46" "0.0692"
Set-TextValue "E46" "  +0.85%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D47" "40.84"
Set-TextValue "E47" "  +0.80%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D48" "342.69"
Set-TextValue "E48" "  +2.57%  "
Set-TextValue "D49" "0.0285"
Set-TextValue "E49" "  +2.34%  "
Set-TextValue "D50" "33.13"
Set-TextValue "E50" "  +8.16%  "
Set-TextValue "E51" "  +6.09%  "
